$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 296.22223
$ws.Range("I127").Value = 257.6154
$ws.Range("J127").Value = 1300
$ws.Range("K127").Value = 772.8462000000001
$ws.Range("L127").Value = 3900
$ws.Range("M127").Value = 4187.1538
$ws.Range("N127").Value = -13820

$ws.Range("H137").Value = 3859.2222
$ws.Range("I137").Value = 4091.625
$ws.Range("K137").Value = 12274.875
$ws.Range("M137").Value = -9724.875

$ws.Range("H138").Value = 1759.7234
$ws.Range("I138").Value = 1534.5416
$ws.Range("J138").Value = 1994.6957
$ws.Range("K138").Value = 4603.6248
$ws.Range("L138").Value = 5984.0871
$ws.Range("M138").Value = 536.3752000000004
$ws.Range("N138").Value = -16264.0871

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8032.0796
$ws.Range("I32").Value = 6821.8037
$ws.Range("K32").Value = 6821.8037
$ws.Range("M32").Value = -6534.8037

$ws.Range("H122").Value = 1214.9584
$ws.Range("I122").Value = 1053.591
$ws.Range("K122").Value = 3160.773
$ws.Range("M122").Value = -710.7729999999997

$ws.Range("H132").Value = 4796.689
$ws.Range("I132").Value = 2822.4211
$ws.Range("J132").Value = 6239.423
$ws.Range("K132").Value = 8467.263300000001
$ws.Range("L132").Value = 18718.269
$ws.Range("M132").Value = -5937.263300000001
$ws.Range("N132").Value = -23778.269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5196
$ws.Range("I22").Value = 5196
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 5196
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -5023
$ws.Range("N22").ClearContents()

$ws.Range("H57").Value = 59600
$ws.Range("J57").Value = 59600
$ws.Range("L57").Value = 59600
$ws.Range("N57").Value = -61040

$ws.Range("H80").Value = 627.9231
$ws.Range("I80").Value = 419.16666
$ws.Range("J80").Value = 806.8570999999999
$ws.Range("K80").Value = 419.16666
$ws.Range("L80").Value = 806.8570999999999
$ws.Range("M80").Value = 578.83334
$ws.Range("N80").Value = -2802.8571

$ws.Range("H83").Value = 627.9231
$ws.Range("I83").Value = 419.16666
$ws.Range("J83").Value = 806.8570999999999
$ws.Range("K83").Value = 2095.8333
$ws.Range("L83").Value = 4034.2855
$ws.Range("M83").Value = 2896.1667
$ws.Range("N83").Value = -14018.2855

$ws.Range("H134").Value = 5802.6665
$ws.Range("I134").Value = 2195.739
$ws.Range("K134").Value = 6587.217000000001
$ws.Range("M134").Value = -4052.217000000001

$ws.Range("H136").Value = 59600
$ws.Range("J136").Value = 59600
$ws.Range("L136").Value = 59600
$ws.Range("N136").Value = -69800

$ws.Range("H137").Value = 50495
$ws.Range("J137").Value = 50495
$ws.Range("L137").Value = 50495
$ws.Range("N137").Value = -60695

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7410135
$ws.Range("I31").Value = 2499.8262
$ws.Range("J31").Value = 15154481
$ws.Range("K31").Value = 2499.8262
$ws.Range("L31").Value = 15154481
$ws.Range("M31").Value = -2204.8262
$ws.Range("N31").Value = -15155071

$ws.Range("H34").Value = 7410135
$ws.Range("I34").Value = 2499.8262
$ws.Range("J34").Value = 15154481
$ws.Range("K34").Value = 2499.8262
$ws.Range("L34").Value = 15154481
$ws.Range("M34").Value = -2297.8262
$ws.Range("N34").Value = -15154885

$ws.Range("H58").Value = 4772.3438
$ws.Range("I58").Value = 2040.4706
$ws.Range("J58").Value = 7868.467
$ws.Range("K58").Value = 2040.4706
$ws.Range("L58").Value = 7868.467
$ws.Range("M58").Value = -1837.4706
$ws.Range("N58").Value = -8274.467000000001

$ws.Range("H86").Value = 5569.381
$ws.Range("I86").Value = 4905.1816
$ws.Range("J86").Value = 6300
$ws.Range("K86").Value = 4905.1816
$ws.Range("L86").Value = 6300
$ws.Range("M86").Value = -3782.1816
$ws.Range("N86").Value = -8546

$ws.Range("H89").Value = 5569.381
$ws.Range("I89").Value = 4905.1816
$ws.Range("J89").Value = 6300
$ws.Range("K89").Value = 24525.908
$ws.Range("L89").Value = 31500
$ws.Range("M89").Value = -18909.908
$ws.Range("N89").Value = -42732

$ws.Range("H132").Value = 2514.6553
$ws.Range("I132").Value = 1681
$ws.Range("J132").Value = 3407.8572
$ws.Range("K132").Value = 5043
$ws.Range("L132").Value = 10223.5716
$ws.Range("M132").Value = -2513
$ws.Range("N132").Value = -15283.5716

$ws.Range("H136").Value = 4772.3438
$ws.Range("I136").Value = 2040.4706
$ws.Range("J136").Value = 7868.467
$ws.Range("K136").Value = 6121.4118
$ws.Range("L136").Value = 23605.401
$ws.Range("M136").Value = -3571.4118
$ws.Range("N136").Value = -28705.401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1320.5652
$ws.Range("J98").Value = 1669.5
$ws.Range("L98").Value = 5008.5
$ws.Range("N98").Value = -8004.5

$ws.Range("H113").Value = 560.2083
$ws.Range("I113").Value = 420.82858
$ws.Range("J113").Value = 935.46155
$ws.Range("K113").Value = 1262.48574
$ws.Range("L113").Value = 2806.38465
$ws.Range("M113").Value = 907.5142599999999
$ws.Range("N113").Value = -7146.38465

$ws.Range("H131").Value = 1141.7632
$ws.Range("J131").Value = 1155.24
$ws.Range("L131").Value = 3465.72
$ws.Range("N131").Value = -13545.72

$ws.Range("H132").Value = 1192.8125
$ws.Range("I132").Value = 867.7143
$ws.Range("J132").Value = 1445.6666
$ws.Range("K132").Value = 7809.428699999999
$ws.Range("L132").Value = 13010.9994
$ws.Range("M132").Value = -5279.428699999999
$ws.Range("N132").Value = -18070.9994

$ws.Range("H134").Value = 3128.2
$ws.Range("I134").Value = 1487.7778
$ws.Range("J134").Value = 5588.8335
$ws.Range("K134").Value = 4463.3334
$ws.Range("L134").Value = 16766.5005
$ws.Range("M134").Value = 606.6665999999996
$ws.Range("N134").Value = -26906.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 20153.846
$ws.Range("J96").Value = 20153.846
$ws.Range("L96").Value = 20153.846
$ws.Range("N96").Value = -25645.846

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 2433.3333
$ws.Range("I39").Value = 2150
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 2150
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -1737
$ws.Range("N39").Value = -3826
